$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1171.4286
$ws.Range("I18").Value = 1240
$ws.Range("K18").Value = 1240
$ws.Range("M18").Value = -956
$ws.Range("H64").Value = 7429.0835
$ws.Range("I64").Value = 5577.5
$ws.Range("J64").Value = 8354.875
$ws.Range("K64").Value = 5577.5
$ws.Range("L64").Value = 8354.875
$ws.Range("M64").Value = -5329.5
$ws.Range("N64").Value = -8850.875
$ws.Range("H67").Value = 7429.0835
$ws.Range("I67").Value = 5577.5
$ws.Range("J67").Value = 8354.875
$ws.Range("K67").Value = 5577.5
$ws.Range("L67").Value = 8354.875
$ws.Range("M67").Value = -4719.5
$ws.Range("N67").Value = -10070.875
$ws.Range("H88").Value = 9580
$ws.Range("I88").Value = 10225
$ws.Range("J88").Value = 7000
$ws.Range("K88").Value = 10225
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = -9819
$ws.Range("N88").Value = -7812
$ws.Range("H91").Value = 9580
$ws.Range("I91").Value = 10225
$ws.Range("J91").Value = 7000
$ws.Range("K91").Value = 10225
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = -8821
$ws.Range("N91").Value = -9808
$ws.Range("H132").Value = 18177.787
$ws.Range("I132").Value = 1779.5769
$ws.Range("K132").Value = 5338.7307
$ws.Range("M132").Value = -2808.7307
$ws.Range("H137").Value = 5721443.5
$ws.Range("J137").Value = 3037.0833
$ws.Range("L137").Value = 9111.249899999999
$ws.Range("N137").Value = -14211.2499
$ws.Range("H138").Value = 9445.896000000001
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 9445.896000000001
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 28337.688
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -38617.688
$ws.Range("H141").Value = 4611.6924
$ws.Range("I141").Value = 3328.9697
$ws.Range("J141").Value = 11666.667
$ws.Range("K141").Value = 9986.909100000001
$ws.Range("L141").Value = 35000.001
$ws.Range("M141").Value = -4806.909100000001
$ws.Range("N141").Value = -45360.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13704600
$ws.Range("I32").Value = 15157827
$ws.Range("K32").Value = 15157827
$ws.Range("M32").Value = -15157540
$ws.Range("H74").Value = 1656.9259
$ws.Range("I74").Value = 1070.5
$ws.Range("J74").Value = 3332.4285
$ws.Range("K74").Value = 1070.5
$ws.Range("L74").Value = 3332.4285
$ws.Range("M74").Value = -196.5
$ws.Range("N74").Value = -5080.4285
$ws.Range("H77").Value = 1656.9259
$ws.Range("I77").Value = 1070.5
$ws.Range("J77").Value = 3332.4285
$ws.Range("K77").Value = 5352.5
$ws.Range("L77").Value = 16662.1425
$ws.Range("M77").Value = -984.5
$ws.Range("N77").Value = -25398.1425
$ws.Range("H132").Value = 2657.2
$ws.Range("I132").Value = 2617.5557
$ws.Range("J132").Value = 3014
$ws.Range("K132").Value = 7852.6671
$ws.Range("L132").Value = 9042
$ws.Range("M132").Value = -5322.6671
$ws.Range("N132").Value = -14102

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 315.33334
$ws.Range("I80").Value = 428.33334
$ws.Range("J80").Value = 301.20834
$ws.Range("K80").Value = 428.33334
$ws.Range("L80").Value = 301.20834
$ws.Range("M80").Value = 569.66666
$ws.Range("N80").Value = -2297.20834
$ws.Range("H83").Value = 315.33334
$ws.Range("I83").Value = 428.33334
$ws.Range("J83").Value = 301.20834
$ws.Range("K83").Value = 2141.6667
$ws.Range("L83").Value = 1506.0417
$ws.Range("M83").Value = 2850.3333
$ws.Range("N83").Value = -11490.0417
$ws.Range("H134").Value = 2961.8147
$ws.Range("I134").Value = 2786.606
$ws.Range("J134").Value = 3237.1428
$ws.Range("K134").Value = 8359.818000000001
$ws.Range("L134").Value = 9711.428400000001
$ws.Range("M134").Value = -5824.818000000001
$ws.Range("N134").Value = -14781.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 461
$ws.Range("I7").Value = 400
$ws.Range("J7").Value = 466.54544
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 466.54544
$ws.Range("M7").Value = -287
$ws.Range("N7").Value = -692.54544
$ws.Range("H16").Value = 10046.538
$ws.Range("I16").Value = 7660.6
$ws.Range("K16").Value = 7660.6
$ws.Range("M16").Value = -7373.6
$ws.Range("H31").Value = 2484.25
$ws.Range("I31").Value = 1998.9736
$ws.Range("J31").Value = 3026.6177
$ws.Range("K31").Value = 1998.9736
$ws.Range("L31").Value = 3026.6177
$ws.Range("M31").Value = -1703.9736
$ws.Range("N31").Value = -3616.6177
$ws.Range("H34").Value = 2484.25
$ws.Range("I34").Value = 1998.9736
$ws.Range("J34").Value = 3026.6177
$ws.Range("K34").Value = 1998.9736
$ws.Range("L34").Value = 3026.6177
$ws.Range("M34").Value = -1796.9736
$ws.Range("N34").Value = -3430.6177
$ws.Range("H41").Value = 4500
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 13166.667
$ws.Range("I50").Value = 13166.667
$ws.Range("K50").Value = 13166.667
$ws.Range("M50").Value = -12541.667
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 59433.25
$ws.Range("J74").Value = 59433.25
$ws.Range("L74").Value = 59433.25
$ws.Range("N74").Value = -61181.25
$ws.Range("H77").Value = 59433.25
$ws.Range("J77").Value = 59433.25
$ws.Range("L77").Value = 178299.75
$ws.Range("N77").Value = -187035.75
$ws.Range("H113").Value = 10046.538
$ws.Range("I113").Value = 7660.6
$ws.Range("K113").Value = 7660.6
$ws.Range("M113").Value = -5490.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 5041.95
$ws.Range("I122").Value = 882.6667
$ws.Range("J122").Value = 8445
$ws.Range("K122").Value = 7944.0003
$ws.Range("L122").Value = 76005
$ws.Range("M122").Value = -5494.0003
$ws.Range("N122").Value = -80905
$ws.Range("H140").Value = 17334812
$ws.Range("I140").Value = 17334812
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 52004436
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -51999256
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 24999.75
$ws.Range("J15").Value = 24999.75
$ws.Range("L15").Value = 24999.75
$ws.Range("N15").Value = -25575.75
$ws.Range("H81").Value = 24999.75
$ws.Range("J81").Value = 24999.75
$ws.Range("L81").Value = 24999.75
$ws.Range("N81").Value = -26995.75
$ws.Range("H84").Value = 24999.75
$ws.Range("J84").Value = 24999.75
$ws.Range("L84").Value = 74999.25
$ws.Range("N84").Value = -84983.25
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 84980
$ws.Range("J110").Value = 84980
$ws.Range("L110").Value = 84980
$ws.Range("N110").Value = -93160
$ws.Range("H132").Value = 4980.2173
$ws.Range("I132").Value = 4691.9697
$ws.Range("J132").Value = 5711.923
$ws.Range("K132").Value = 14075.9091
$ws.Range("L132").Value = 17135.769
$ws.Range("M132").Value = -11545.9091
$ws.Range("N132").Value = -22195.769
$ws.Range("H136").Value = 4238.15
$ws.Range("I136").Value = 2619.7
$ws.Range("J136").Value = 5856.6
$ws.Range("K136").Value = 7859.099999999999
$ws.Range("L136").Value = 17569.8
$ws.Range("M136").Value = -5309.099999999999
$ws.Range("N136").Value = -22669.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 55499.832
$ws.Range("I75").Value = 39999
$ws.Range("J75").Value = 58600
$ws.Range("K75").Value = 39999
$ws.Range("L75").Value = 58600
$ws.Range("M75").Value = -39063
$ws.Range("N75").Value = -60472
$ws.Range("H78").Value = 55499.832
$ws.Range("I78").Value = 39999
$ws.Range("J78").Value = 58600
$ws.Range("K78").Value = 119997
$ws.Range("L78").Value = 175800
$ws.Range("M78").Value = -115317
$ws.Range("N78").Value = -185160
$ws.Range("H80").Value = 109499.5
$ws.Range("J80").Value = 109499.5
$ws.Range("L80").Value = 109499.5
$ws.Range("N80").Value = -111495.5
$ws.Range("H83").Value = 109499.5
$ws.Range("J83").Value = 109499.5
$ws.Range("L83").Value = 328498.5
$ws.Range("N83").Value = -338482.5
$ws.Range("H132").Value = 3608.12
$ws.Range("I132").Value = 3238.6316
$ws.Range("J132").Value = 4778.1665
$ws.Range("K132").Value = 9715.8948
$ws.Range("L132").Value = 14334.4995
$ws.Range("M132").Value = -7185.8948
$ws.Range("N132").Value = -19394.4995
$ws.Range("H136").Value = 7679.72
$ws.Range("I136").Value = 7661.5713
$ws.Range("K136").Value = 22984.7139
$ws.Range("M136").Value = -20434.7139
